# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a different (Strike#) quantity; this
# recomputes/overwrites it with the correct strikeout (K) counts for each
# start, row by row (rows 2-35; row 36 already held the correct value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 6
    6  = 3
    7  = 5
    8  = 5
    9  = 4
    10 = 1
    11 = 4
    12 = 5
    13 = 3
    14 = 7
    15 = 4
    16 = 6
    17 = 7
    18 = 6
    19 = 9
    20 = 5
    21 = 9
    22 = 4
    23 = 3
    24 = 2
    25 = 7
    26 = 5
    27 = 3
    28 = 4
    29 = 9
    30 = 5
    31 = 4
    32 = 8
    33 = 4
    34 = 8
    35 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
